# Add 2022-Q4 data
# ------------------------------------------------------------------
# Strategy:
#  - "总计" sheet (summary): insert a new row for 2022-Q4 data, shifting
#    the existing 2022-Q3 / 2022-Q2 rows down by one.
#  - The worksheet that used to be named "2022-Q3" is repurposed in
#    place to hold the brand-new 2022-Q4 fund data (renamed accordingly).
#  - A duplicate of the original "2022-Q3" sheet (with its original
#    fund data untouched) is inserted right after it and renamed back
#    to "2022-Q3", preserving that historical data.
#  - The "2022-Q2" sheet is left completely untouched.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "总计" (summary) sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Insert a fresh row 3 - this pushes the current row3 (2022-Q2) down to
# row4, and leaves row2 (currently 2022-Q3) in place for now.
$summary.Rows.Item(3).Insert()

# Give the newly inserted A3 the same look (bold / centered / bordered)
# as the other index cells in column A (A2, A4, ...) by copying the
# formatting from A2.
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

# Row 3 now becomes the historical "2022-Q3" entry (what row2 used to
# hold), row 2 becomes the new "2022-Q4" entry, row 4 keeps "2022-Q2".
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.63

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.05

$summary.Range("A4").Value = 2

# ---------------------------------------------------------------
# 2) Duplicate the current "2022-Q3" sheet so its original data is
#    preserved under the same name after we repurpose the original.
# ---------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($null, $q3)
$q3Dup = $wb.Worksheets.Item("2022-Q3 (2)")
$q3Dup.Name = "2022-Q3-new"

# ---------------------------------------------------------------
# 3) Repurpose the original "2022-Q3" sheet to hold the new
#    "2022-Q4" fund data.
# ---------------------------------------------------------------
$q3.Name = "2022-Q4"

# Columns B-G hold text data (fund code / name / figures kept as
# strings, exactly like the rest of the workbook) - force a text
# number format first so Excel does not auto-convert numeric-looking
# strings (e.g. "009246") into numbers.
$q3.Range("B2:G4").NumberFormat = "@"

$q3.Range("B2").Value = "009246"
$q3.Range("C2").Value = "大摩ESG量化混合"
$q3.Range("D2").Value = "2.39"
$q3.Range("E2").Value = "92.14"
$q3.Range("F2").Value = "1.17"
$q3.Range("G2").Value = "0.0280"
$q3.Range("H2").Value = 8

$q3.Range("B3").Value = "620002"
$q3.Range("C3").Value = "金元顺安成长动力混合"
$q3.Range("D3").Value = "0.35"
$q3.Range("E3").Value = "72.02"
$q3.Range("F3").Value = "3.14"
$q3.Range("G3").Value = "0.0110"
$q3.Range("H3").Value = 8

$q3.Range("B4").Value = "165531"
$q3.Range("C4").Value = "信诚多策略灵活配置混合（LOF）"
$q3.Range("D4").Value = "0.89"
$q3.Range("E4").Value = "72.25"
$q3.Range("F4").Value = "1.04"
$q3.Range("G4").Value = "0.0093"
$q3.Range("H4").Value = 7

# ---------------------------------------------------------------
# 4) Rename the duplicated sheet back to "2022-Q3" - it already has
#    the original (now-historical) fund data.
# ---------------------------------------------------------------
$q3Dup.Name = "2022-Q3"

# ---------------------------------------------------------------
# 5) Restore the original active tab (2022-Q2, the last sheet).
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Activate()
